# Update crypto price/volume table with latest values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.NumberFormat = "General"
}

Set-TextValue "D2" "59.516.99"
Set-TextValue "E2" "  -0.69%  "
Set-TextValue "D3" "2.650.17"
Set-TextValue "E3" "  -0.45%  "
Set-TextValue "E4" "  -0.18%  "
Set-TextValue "D5" "518.45"
Set-TextValue "E5" "  -0.68%  "
Set-TextValue "D6" "146.58"
Set-TextValue "E6" "  -1.71%  "
Set-TextValue "E7" "  +0.26%  "
Set-TextValue "D8" "0.573"
Set-TextValue "E8" "  -0.28%  "
Set-TextValue "D9" "2.660.05"
Set-TextValue "E9" "  -1.08%  "
Set-TextValue "D10" "6.32"
Set-TextValue "E10" "  -3.27%  "
Set-TextValue "D11" "0.106"
Set-TextValue "E11" "  -0.86%  "
Set-TextValue "D12" "0.338"
Set-TextValue "E12" "  -1.34%  "
Set-TextValue "E13" "  +0.51%  "
Set-TextValue "D14" "3.112.51"
Set-TextValue "E14" "  -0.42%  "
Set-TextValue "D15" "59.481.19"
Set-TextValue "E15" "  -0.50%  "
Set-TextValue "D16" "20.99"
Set-TextValue "E16" "  -2.70%  "
Set-TextValue "E17" "  -1.37%  "
Set-TextValue "D18" "2.651.33"
Set-TextValue "E18" "  -0.90%  "
Set-TextValue "D19" "350.54"
Set-TextValue "E19" "  +0.56%  "
Set-TextValue "D20" "4.51"
Set-TextValue "E20" "  -2.85%  "
Set-TextValue "D21" "10.37"
Set-TextValue "E21" "  -2.75%  "
Set-TextValue "D22" "6.22"
Set-TextValue "E22" "  -0.86%  "
Set-TextValue "D23" "1.00"
Set-TextValue "E23" "  +0.18%  "
Set-TextValue "D24" "62.05"
Set-TextValue "E24" "  +1.36%  "
Set-TextValue "D25" "0.417"
Set-TextValue "E25" "  -2.77%  "
Set-TextValue "E26" "  +1.72%  "
Set-TextValue "E27" "  +0.08%  "
Set-TextValue "D28" "0.0₃0809"
Set-TextValue "E28" "  -3.40%  "
Set-TextValue "D29" "7.17"
Set-TextValue "E29" "  -1.13%  "
Set-TextValue "E30" "  +0.21%  "
Set-TextValue "D31" "6.49"
Set-TextValue "E31" "  -2.10%  "
Set-TextValue "E32" "  -1.11%  "
Set-TextValue "D33" "18.96"
Set-TextValue "E33" "  -1.02%  "
Set-TextValue "D34" "149.22"
Set-TextValue "E34" "  -0.01%  "
Set-TextValue "D35" "4.07"
Set-TextValue "E35" "  +0.07%  "
Set-TextValue "D36" "0.941"
Set-TextValue "E36" "  -13.34%  "
Set-TextValue "E37" "  +0.27%  "
Set-TextValue "D38" "0.869"
Set-TextValue "E38" "  -1.40%  "
Set-TextValue "D39" "36.70"
Set-TextValue "E39" "  +0.07%  "
Set-TextValue "D40" "1.47"
Set-TextValue "E40" "  +1.60%  "
Set-TextValue "D41" "3.69"
Set-TextValue "E41" "  -1.36%  "
Set-TextValue "D42" "279.08"
Set-TextValue "E42" "  -4.17%  "
Set-TextValue "D43" "0.0994"
Set-TextValue "E43" "  -0.61%  "
Set-TextValue "E44" "  +0.51%  "
Set-TextValue "D45" "0.605"
Set-TextValue "E45" "  -4.13%  "
Set-TextValue "D46" "19.74"
Set-TextValue "E46" "  -0.35%  "
Set-TextValue "D47" "2.101.57"
Set-TextValue "E47" "  +5.46%  "
Set-TextValue "D48" "0.0530"
Set-TextValue "E48" "  -3.87%  "
Set-TextValue "E49" "  -0.86%  "
Set-TextValue "E50" "  -1.79%  "
Set-TextValue "E51" "  +0.46%  "
